# HashMap to read testdata
# Updates the "smoke" sheet test data so the user/password/leave fields are
# written in the "key:=value" HashMap style the automation framework reads.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("smoke")
$ws.Activate()

# ---------------------------------------------------------------------
# 1. Data cells: prefix values with their HashMap keys.
# ---------------------------------------------------------------------

# Row 2 (tc001_VerifyLoginAuthentication)
$ws.Range("C2").Value = "user:=user02"
$ws.Range("D2").Value = "pwd:=TM1234"

# Row 3 (tc002_ApplyLeave)
$ws.Range("C3").Value = "user:=user03"
$ws.Range("D3").Value = "pwd:=TM1235"
$ws.Range("E3").Value = "leaveType:=Annual Leave"
$ws.Range("F3").Value = "fromDate:=25-07-2020"
$ws.Range("G3").Value = "toDate:=26-07-2020"
$ws.Range("H3").Value = "comment:=attend function"

# Row 4 (tc003_CancelLeave)
$ws.Range("C4").Value = "user:=user04"
$ws.Range("D4").Value = "pwd:=TM1236"

# ---------------------------------------------------------------------
# 2. Shrink the header/data font from 18pt to 10pt.
# ---------------------------------------------------------------------

$ws.Range("A1:H4").Font.Size = 10

# Row heights were manually set (23.4) to fit the old 18pt font; now that
# the font is smaller, auto-fit each row so Excel drops the explicit
# height and falls back to the sheet default.
for ($r = 1; $r -le 4; $r++) {
    $ws.Rows.Item($r).AutoFit()
}

# Sheet default row height shrinks to match the smaller default font.
$ws.StandardHeight = 13.8

# ---------------------------------------------------------------------
# 3. Column widths - re-fit now that the font/content changed.
# ---------------------------------------------------------------------

$ws.Columns.Item(1).ColumnWidth = 4.33203125
$ws.Columns.Item(2).ColumnWidth = 27
$ws.Columns.Item(3).ColumnWidth = 11.109375
$ws.Columns.Item(4).ColumnWidth = 12.33203125
$ws.Columns.Item(5).ColumnWidth = 19.6640625
$ws.Columns.Item(6).ColumnWidth = 14.109375
$ws.Columns.Item(7).ColumnWidth = 12.6640625
$ws.Columns.Item(8).ColumnWidth = 22.21875

# ---------------------------------------------------------------------
# 4. Selection moves from C7 to the whole of row 3.
# ---------------------------------------------------------------------

$ws.Rows.Item(3).Select()
